$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    9  = -21.782
    13 = -22.08
    16 = -22.002
    18 = -21.982
    20 = -20.623
    26 = -21.647
    27 = -21.839
    29 = -21.201
    35 = -20.123
    36 = -20.391
    45 = -21.469
    55 = -22.124
    57 = -22.232
    69 = -21.52
    76 = -20.047
    78 = -19.854
    82 = -21.962
    83 = -21.891
    93 = -21.508
    97 = -21.62
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 1).Value = $updates[$row]
}
